$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Cruz Laminates & Countertops", "5207 Harrisburg Blvd, Houston, TX 77011"),
    @("Sweeney Marble Co", "2401 Polk St, Houston, TX 77003"),
    @("Victor's Granite Designers", "3822 Ranch St, Houston, TX 77026"),
    @("Texas Custom Marble & Granite", "10835 Maple Leaf St, Houston, TX 77016"),
    @("Floor & Decor", "4330 Dacoma St, Houston, TX 77092")
)

$startRow = 474
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
